$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 4026.8667
$ws.Range("I74").Value = 3925.375
$ws.Range("K74").Value = 3925.375
$ws.Range("M74").Value = -2989.375

$ws.Range("H77").Value = 4026.8667
$ws.Range("I77").Value = 3925.375
$ws.Range("K77").Value = 19626.875
$ws.Range("M77").Value = -14946.875

$ws.Range("H129").Value = 1117.0488
$ws.Range("I129").Value = 632.3333
$ws.Range("J129").Value = 1155.3158
$ws.Range("K129").Value = 1896.9999
$ws.Range("L129").Value = 3465.9474
$ws.Range("M129").Value = 3103.0001
$ws.Range("N129").Value = -13465.9474

$ws.Range("H132").Value = 29936.361
$ws.Range("I132").Value = 30763.115
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 92289.345
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -89759.345
$ws.Range("N132").Value = -8060

$ws.Range("H135").Value = 413.48718
$ws.Range("I135").Value = 394.33334
$ws.Range("J135").Value = 643.3333
$ws.Range("K135").Value = 3549.00006
$ws.Range("L135").Value = 5789.9997
$ws.Range("M135").Value = -1014.00006
$ws.Range("N135").Value = -10859.9997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4539.1235
$ws.Range("I32").Value = 3297.13
$ws.Range("J32").Value = 28447.5
$ws.Range("K32").Value = 3297.13
$ws.Range("L32").Value = 28447.5
$ws.Range("M32").Value = -3010.13
$ws.Range("N32").Value = -29021.5

$ws.Range("H57").Value = 3000
$ws.Range("I57").Value = 3000
$ws.Range("K57").Value = 3000
$ws.Range("M57").Value = -2516

$ws.Range("H61").Value = 1396.9318
$ws.Range("I61").Value = 865.5925999999999
$ws.Range("J61").Value = 2240.8235
$ws.Range("K61").Value = 865.5925999999999
$ws.Range("L61").Value = 2240.8235
$ws.Range("M61").Value = -653.5925999999999
$ws.Range("N61").Value = -2664.8235

$ws.Range("H97").Value = 1290.7142
$ws.Range("I97").Value = 1046
$ws.Range("J97").Value = 1617
$ws.Range("K97").Value = 1046
$ws.Range("L97").Value = 1617
$ws.Range("M97").Value = -550
$ws.Range("N97").Value = -2609

$ws.Range("H132").Value = 1714.0526
$ws.Range("I132").Value = 875.7879
$ws.Range("J132").Value = 2866.6667
$ws.Range("K132").Value = 2627.3637
$ws.Range("L132").Value = 8600.000100000001
$ws.Range("M132").Value = -97.36369999999988
$ws.Range("N132").Value = -13660.0001

$ws.Range("H136").Value = 1396.9318
$ws.Range("I136").Value = 865.5925999999999
$ws.Range("J136").Value = 2240.8235
$ws.Range("K136").Value = 2596.7778
$ws.Range("L136").Value = 6722.470499999999
$ws.Range("M136").Value = -46.77779999999984
$ws.Range("N136").Value = -11822.4705

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H30").Value = 50000
$ws.Range("J30").Value = 50000
$ws.Range("L30").Value = 50000
$ws.Range("N30").Value = -50250

$ws.Range("H94").Value = 495.5
$ws.Range("I94").Value = 495.5
$ws.Range("K94").Value = 495.5
$ws.Range("M94").Value = -44.5

$ws.Range("H107").Value = 577.3
$ws.Range("I107").Value = 450.17648
$ws.Range("J107").Value = 1297.6666
$ws.Range("K107").Value = 450.17648
$ws.Range("L107").Value = 1297.6666
$ws.Range("M107").Value = 1469.82352
$ws.Range("N107").Value = -5137.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H76").Value = 1440
$ws.Range("I76").Value = 1440
$ws.Range("K76").Value = 1440
$ws.Range("M76").Value = -1125

$ws.Range("H79").Value = 1440
$ws.Range("I79").Value = 1440
$ws.Range("K79").Value = 1440
$ws.Range("M79").Value = -348

$ws.Range("H107").Value = 823.625
$ws.Range("I107").Value = 549.25
$ws.Range("J107").Value = 1372.375
$ws.Range("K107").Value = 549.25
$ws.Range("L107").Value = 1372.375
$ws.Range("M107").Value = 1370.75
$ws.Range("N107").Value = -5212.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1363.4286
$ws.Range("I113").Value = 1659.7778
$ws.Range("J113").Value = 830
$ws.Range("K113").Value = 4979.3334
$ws.Range("L113").Value = 2490
$ws.Range("M113").Value = -2809.3334
$ws.Range("N113").Value = -6830

$ws.Range("H129").Value = 2387.7058
$ws.Range("I129").Value = 917.8570999999999
$ws.Range("J129").Value = 3416.6
$ws.Range("K129").Value = 2753.5713
$ws.Range("L129").Value = 10249.8
$ws.Range("M129").Value = 2246.4287
$ws.Range("N129").Value = -20249.8

$ws.Range("H131").Value = 4612.1724
$ws.Range("I131").Value = 661.1111
$ws.Range("J131").Value = 6390.15
$ws.Range("K131").Value = 1983.3333
$ws.Range("L131").Value = 19170.45
$ws.Range("M131").Value = 3056.6667
$ws.Range("N131").Value = -29250.45

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1610
$ws.Range("I107").Value = 2200
$ws.Range("J107").Value = 135
$ws.Range("K107").Value = 2200
$ws.Range("L107").Value = 135
$ws.Range("M107").Value = -280
$ws.Range("N107").Value = -3975

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 898.2
$ws.Range("I22").Value = 999.6667
$ws.Range("J22").Value = 746
$ws.Range("K22").Value = 999.6667
$ws.Range("L22").Value = 746
$ws.Range("M22").Value = -704.6667
$ws.Range("N22").Value = -1336

$ws.Range("H27").Value = 898.2
$ws.Range("I27").Value = 999.6667
$ws.Range("J27").Value = 746
$ws.Range("K27").Value = 999.6667
$ws.Range("L27").Value = 746
$ws.Range("M27").Value = -892.6667
$ws.Range("N27").Value = -960

$ws.Range("H50").Value = 23500
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 23500
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 23500
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -24774

$ws.Range("H55").Value = 478
$ws.Range("I55").Value = 273.66666
$ws.Range("J55").Value = 682.3333
$ws.Range("K55").Value = 273.66666
$ws.Range("L55").Value = 682.3333
$ws.Range("M55").Value = -100.66666
$ws.Range("N55").Value = -1028.3333

$ws.Range("H122").Value = 9263034
$ws.Range("I122").Value = 12347636
$ws.Range("J122").Value = 9226.666999999999
$ws.Range("K122").Value = 37042908
$ws.Range("L122").Value = 27680.001
$ws.Range("M122").Value = -37040458
$ws.Range("N122").Value = -32580.001

$ws.Range("H132").Value = 7272.6294
$ws.Range("I132").Value = 7308.079
$ws.Range("J132").Value = 7188.4375
$ws.Range("K132").Value = 21924.237
$ws.Range("L132").Value = 21565.3125
$ws.Range("M132").Value = -19394.237
$ws.Range("N132").Value = -26625.3125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 7001
$ws.Range("J14").Value = 7001
$ws.Range("L14").Value = 7001
$ws.Range("N14").Value = -7337

$ws.Range("H81").Value = 1616.7858
$ws.Range("I81").Value = 967
$ws.Range("J81").Value = 1794
$ws.Range("K81").Value = 1934
$ws.Range("L81").Value = 3588
$ws.Range("M81").Value = -873
$ws.Range("N81").Value = -5710

$ws.Range("H84").Value = 1616.7858
$ws.Range("I84").Value = 967
$ws.Range("J84").Value = 1794
$ws.Range("K84").Value = 9670
$ws.Range("L84").Value = 17940
$ws.Range("M84").Value = -4366
$ws.Range("N84").Value = -28548

$ws.Range("H131").Value = 29997.5
$ws.Range("J131").Value = 29997.5
$ws.Range("L131").Value = 29997.5
$ws.Range("N131").Value = -40077.5

$ws.Range("H132").Value = 1682.9814
$ws.Range("I132").Value = 1513.579
$ws.Range("J132").Value = 2085.3125
$ws.Range("K132").Value = 4540.737
$ws.Range("L132").Value = 6255.9375
$ws.Range("M132").Value = -2010.737
$ws.Range("N132").Value = -11315.9375
